$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 54, shifting existing rows 54..111 down to 55..112
$newRow = $ws.Rows.Item(54)
$newRow.Insert()

# Fill in the constant columns (same for every record in this sheet)
$ws.Cells.Item(54, 1).Value = 11
$ws.Cells.Item(54, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(54, 3).Value = "Bíobío"
$ws.Cells.Item(54, 4).Value = 44638
$ws.Cells.Item(54, 5).Value = 8
$ws.Cells.Item(54, 6).Value = 100112043
$ws.Cells.Item(54, 7).Value = "Pepino ensalada"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 270
$ws.Cells.Item(54, 11).Value = 18000
$ws.Cells.Item(54, 12).Value = 19000
$ws.Cells.Item(54, 13).Value = 18556
$ws.Cells.Item(54, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(54, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(54, 16).Value = 309
$ws.Cells.Item(54, 17).Value = 60
$ws.Cells.Item(54, 18).Value = "Hortaliza"

# Apply the same date style (numFmt) used on the other D-column cells to the new D54
$ws.Cells.Item(54, 4).NumberFormat = $ws.Cells.Item(53, 4).NumberFormat
